# Usability2 sheet: simplify action/query syntax (commas -> semicolons, drop
# redundant parentheses) and refresh the row heights / view state that Excel
# recalculated as a consequence of the shorter text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shared-string text simplifications
# ---------------------------------------------------------------------------

# "Add to cart" action (E3)
$ws.Range("E3").Value() = "CART(); IF [time] == few_days THEN LOAD(Explanation); JUMP(Safety) ELSE JUMP(Explanation)"

# "IF ... THEN ... ELSE JUMP(Explanation)" action, repeated across many rows
$ifAction = "IF [time] == few_days THEN LOAD(Explanation); JUMP(Safety) ELSE JUMP(Explanation)"
foreach ($r in 17..26) {
    $ws.Range("E$r").Value() = $ifAction
}
foreach ($r in 43,44,47,48) {
    $ws.Range("E$r").Value() = $ifAction
}

# Per-medform SELECT queries (F17:F26)
$ws.Range("F17").Value() = "SELECT * FROM Products WHERE MedFormID == 'Bta' AND APIID == [api_calmer]"
$ws.Range("F18").Value() = "SELECT * FROM Products WHERE MedFormID == 'Eli' AND APIID == [api_calmer]"
$ws.Range("F19").Value() = "SELECT * FROM Products WHERE MedFormID == 'DiGra' AND APIID == [api_calmer]"
$ws.Range("F20").Value() = "SELECT * FROM Products WHERE MedFormID == 'Inh' AND APIID == [api_calmer]"
$ws.Range("F21").Value() = "SELECT * FROM Products WHERE MedFormID == 'Kap' AND APIID == [api_calmer]"
$ws.Range("F22").Value() = "SELECT * FROM Products WHERE MedFormID == 'Li' AND APIID == [api_calmer]"
$ws.Range("F23").Value() = "SELECT * FROM Products WHERE MedFormID == 'Lt' AND APIID == [api_calmer]"
$ws.Range("F24").Value() = "SELECT * FROM Products WHERE MedFormID == 'Sa' AND APIID == [api_calmer]"
$ws.Range("F25").Value() = "SELECT * FROM Products WHERE MedFormID == 'Tab' AND APIID == [api_calmer]"
$ws.Range("F26").Value() = "SELECT * FROM Products WHERE MedFormID == 'Tro' AND APIID == [api_calmer]"

# match_decide CODE action (E42)
$ws.Range("E42").Value() = "{med_forms, matched} = MATCH([api_calmer], [water], [swallow], [transport], [fly], [single]); SAVE(med_forms); SAVE(matched); IF ROWS([matched]) == 0 THEN GO(no_match) ELSE GO(match)"

# same_med_forms rich-text action (F47) - keep the two runs / their formatting,
# only shrink the literal text of each run in place.
$cellF47 = $ws.Range("F47")
$run1 = $cellF47.Characters(1, 77)
$run1.Text = "category = TO_TEXT(SELECT CategoryChoice FROM API WHERE ID == [api_calmer]); "
$run2 = $cellF47.Characters(78, 122)
$run2.Text = "SELECT * FROM Products p JOIN API a ON p.APIID == a.ID WHERE a.CategoryChoice == [category] AND MedFormID IN [med_forms]"
$run2color = $cellF47.Characters(78, 120)
$run2color.Font.Color = 0

# ---------------------------------------------------------------------------
# 2. View state (scroll position / active cell) refreshed by Excel
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 3
$ws.Range("F53").Select()

# ---------------------------------------------------------------------------
# 3. Default column width tweak
# ---------------------------------------------------------------------------
$ws.StandardWidth = 11.39453125

# ---------------------------------------------------------------------------
# 4. Row height adjustments (Excel recalculated wrapped-text row heights after
#    the shorter strings were entered)
# ---------------------------------------------------------------------------
foreach ($r in 3,17,18,19,20,21,22,23,24,25,26,43,44,47,48) {
    $ws.Rows.Item($r).RowHeight = 41.75
}
$ws.Rows.Item(42).RowHeight = 68.65
# row 3 is 41.75 too (see above loop); row 42 handled separately since its
# target height differs (68.65).

# ---------------------------------------------------------------------------
# 5. E48 style: align with the plain "Normal" style used by the sibling rows
#    (previously it referenced a near-duplicate style index).
# ---------------------------------------------------------------------------
$ws.Range("E48").Style = $ws.Range("E17").Style
